$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4370.1763
$ws.Range("I43").Value = 3710.5557
$ws.Range("J43").Value = 5112.25
$ws.Range("K43").Value = 3710.5557
$ws.Range("L43").Value = 5112.25
$ws.Range("M43").Value = -3641.5557
$ws.Range("N43").Value = -5250.25
$ws.Range("H134").Value = 39818.184
$ws.Range("J134").Value = 39818.184
$ws.Range("L134").Value = 39818.184
$ws.Range("N134").Value = -49958.184
$ws.Range("H137").Value = 5731.5
$ws.Range("I137").Value = 6326.263
$ws.Range("J137").Value = 4704.1816
$ws.Range("K137").Value = 18978.789
$ws.Range("L137").Value = 14112.5448
$ws.Range("M137").Value = -16428.789
$ws.Range("N137").Value = -19212.5448
$ws.Range("H138").Value = 3654.6897
$ws.Range("I138").Value = 1870.7142
$ws.Range("K138").Value = 5612.142599999999
$ws.Range("M138").Value = -472.1425999999992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1229.5952
$ws.Range("I74").Value = 1160.4242
$ws.Range("J74").Value = 1483.2222
$ws.Range("K74").Value = 1160.4242
$ws.Range("L74").Value = 1483.2222
$ws.Range("M74").Value = -286.4241999999999
$ws.Range("N74").Value = -3231.2222
$ws.Range("H77").Value = 1229.5952
$ws.Range("I77").Value = 1160.4242
$ws.Range("J77").Value = 1483.2222
$ws.Range("K77").Value = 5802.120999999999
$ws.Range("L77").Value = 7416.111
$ws.Range("M77").Value = -1434.120999999999
$ws.Range("N77").Value = -16152.111
$ws.Range("H122").Value = 3797.3635
$ws.Range("I122").Value = 3411.5925
$ws.Range("J122").Value = 5533.3335
$ws.Range("K122").Value = 10234.7775
$ws.Range("L122").Value = 16600.0005
$ws.Range("M122").Value = -7784.7775
$ws.Range("N122").Value = -21500.0005
$ws.Range("H132").Value = 21278436
$ws.Range("I132").Value = 25642564
$ws.Range("K132").Value = 76927692
$ws.Range("M132").Value = -76925162

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 15626350
$ws.Range("I86").Value = 19232024
$ws.Range("J86").Value = 1762.1666
$ws.Range("K86").Value = 19232024
$ws.Range("L86").Value = 1762.1666
$ws.Range("M86").Value = -19230901
$ws.Range("N86").Value = -4008.1666
$ws.Range("H89").Value = 15626350
$ws.Range("I89").Value = 19232024
$ws.Range("J89").Value = 1762.1666
$ws.Range("K89").Value = 96160120
$ws.Range("L89").Value = 8810.833000000001
$ws.Range("M89").Value = -96154504
$ws.Range("N89").Value = -20042.833
$ws.Range("H94").Value = 2065.1738
$ws.Range("I94").Value = 791.4666999999999
$ws.Range("K94").Value = 791.4666999999999
$ws.Range("M94").Value = -340.4666999999999
$ws.Range("H134").Value = 2012.4546
$ws.Range("I134").Value = 1126.3334
$ws.Range("K134").Value = 3379.0002
$ws.Range("M134").Value = -844.0001999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I14").Value = 650
$ws.Range("J14").Value = 2500
$ws.Range("K14").Value = 650
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = -480
$ws.Range("N14").Value = -2840
$ws.Range("H58").Value = 1617.0588
$ws.Range("I58").Value = 1462
$ws.Range("K58").Value = 1462
$ws.Range("M58").Value = -1259
$ws.Range("H132").Value = 4122.647
$ws.Range("I132").Value = 2997.6365
$ws.Range("K132").Value = 8992.9095
$ws.Range("M132").Value = -6462.9095
$ws.Range("H134").Value = 3314.5715
$ws.Range("I134").Value = 2918.0715
$ws.Range("J134").Value = 4107.5713
$ws.Range("K134").Value = 8754.2145
$ws.Range("L134").Value = 12322.7139
$ws.Range("M134").Value = -6219.2145
$ws.Range("N134").Value = -17392.7139
$ws.Range("H136").Value = 1617.0588
$ws.Range("I136").Value = 1462
$ws.Range("K136").Value = 4386
$ws.Range("M136").Value = -1836

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 11111464
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 12500387
$ws.Range("K2").Value = 80
$ws.Range("L2").Value = 12500387
$ws.Range("M2").Value = 33
$ws.Range("N2").Value = -12500613
$ws.Range("H7").Value = 5037500
$ws.Range("J7").Value = 5050000
$ws.Range("L7").Value = 5050000
$ws.Range("N7").Value = -5050224
$ws.Range("H8").Value = 5037500
$ws.Range("J8").Value = 5050000
$ws.Range("L8").Value = 5050000
$ws.Range("N8").Value = -5050278
$ws.Range("H11").Value = 2022039.2
$ws.Range("I11").Value = 1293981.5
$ws.Range("J11").Value = 3023118.5
$ws.Range("K11").Value = 1293981.5
$ws.Range("L11").Value = 3023118.5
$ws.Range("M11").Value = -1293842.5
$ws.Range("N11").Value = -3023396.5
$ws.Range("H14").Value = 864090.3
$ws.Range("I14").Value = 1876225
$ws.Range("J14").Value = 189333.83
$ws.Range("K14").Value = 1876225
$ws.Range("L14").Value = 189333.83
$ws.Range("M14").Value = -1876057
$ws.Range("N14").Value = -189669.83
$ws.Range("H126").Value = 712.1667
$ws.Range("I126").Value = 769.8
$ws.Range("K126").Value = 2309.4
$ws.Range("M126").Value = 160.6000000000004
$ws.Range("H132").Value = 6343.4146
$ws.Range("I132").Value = 6140.552
$ws.Range("J132").Value = 6833.6665
$ws.Range("K132").Value = 18421.656
$ws.Range("L132").Value = 20500.9995
$ws.Range("M132").Value = -15891.656
$ws.Range("N132").Value = -25560.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 333.64865
$ws.Range("I55").Value = 353.89655
$ws.Range("J55").Value = 260.25
$ws.Range("K55").Value = 353.89655
$ws.Range("L55").Value = 260.25
$ws.Range("M55").Value = -180.89655
$ws.Range("N55").Value = -606.25
$ws.Range("H122").Value = 41671150
$ws.Range("I122").Value = 90911520
$ws.Range("K122").Value = 272734560
$ws.Range("M122").Value = -272732110

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 127
$ws.Range("I7").Value = 127
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 127
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -14
$ws.Range("N7").ClearContents()
$ws.Range("H113").Value = 554.1429000000001
$ws.Range("I113").Value = 456
$ws.Range("K113").Value = 1368
$ws.Range("M113").Value = 802
$ws.Range("H126").Value = 6411.385
$ws.Range("I126").Value = 9335.429
$ws.Range("K126").Value = 28006.287
$ws.Range("M126").Value = -25536.287
